$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on D2:E51 so numeric-looking strings
# (e.g. "584.49", "0.995") are written back as text, matching the
# original inline-string cell content instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.364.80"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "2.655.75"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "584.49"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "144.57"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").Value = "0.382"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").Value = "0.157"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "3.127.70"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "26.15"
$ws.Range("E14").Value = "  +7.08%  "
$ws.Range("D15").Value = "61.280.87"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("D17").Value = "2.666.22"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "11.68"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "4.79"
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("D20").Value = "355.89"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "0.524"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "64.64"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("E26").Value = "  +4.97%  "
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "2.00"
$ws.Range("E28").Value = "  +7.92%  "
$ws.Range("D29").Value = "0.0₃0825"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  +8.70%  "
$ws.Range("D31").Value = "169.48"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "20.16"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("E34").Value = "  +14.57%  "
$ws.Range("D35").Value = "4.68"
$ws.Range("E35").Value = "  +8.08%  "
$ws.Range("E36").Value = "  +7.70%  "
$ws.Range("E37").Value = "  +19.09%  "
$ws.Range("D38").Value = "1.72"
$ws.Range("E38").Value = "  +4.91%  "
$ws.Range("D39").Value = "345.86"
$ws.Range("E39").Value = "  +9.87%  "
$ws.Range("E40").Value = "  +6.13%  "
$ws.Range("D41").Value = "38.52"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +6.37%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "20.81"
$ws.Range("E43").Value = "  +4.69%  "
$ws.Range("D44").Value = "21.28"
$ws.Range("E44").Value = "  +5.88%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0579"
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.630"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0254"
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "135.71"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "2.101.55"
$ws.Range("E51").Value = "  +3.40%  "

# Restore the default (unstyled) cell style so no stray formatting
# is introduced versus the original workbook.
$ws.Range("D2:E51").Style = "Normal"
